$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 30-32: Approved/Rejected column (I) flips from "Approved" to "Rejected",
# and a "Nil" reason is added in the ReasonToReject column (J).
$ws.Range("I30").Value = "Rejected"
$ws.Range("J30").Value = "Nil"

$ws.Range("I31").Value = "Rejected"
$ws.Range("J31").Value = "Nil"

$ws.Range("I32").Value = "Rejected"
$ws.Range("J32").Value = "Nil"

# Row 36: Approved/Rejected column (I) flips from "Rejected" to "Approved",
# and its ReasonToReject value (J) is cleared out.
$ws.Range("I36").Value = "Approved"
$ws.Range("J36").ClearContents()

# Update the active selection to reflect where the user ended up (H37).
$ws.Range("H37").Select()
